$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.649.64'
$ws.Range("E2").Value = '  +2.97%  '
$ws.Range("D3").Value = '2.196.83'
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'259.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.88%  '
$ws.Range("D6").Value = "'82.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +12.35%  '
$ws.Range("E7").Value = '  +1.79%  '
$ws.Range("E9").Value = '  +2.00%  '
$ws.Range("D10").Value = "'43.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.63%  '
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("D12").Value = "'6.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.64%  '
$ws.Range("E13").Value = '  +1.98%  '
$ws.Range("D14").Value = '2.523.91'
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("D15").Value = "'14.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").Value = '2.205.05'
$ws.Range("E16").Value = '  +1.51%  '
$ws.Range("D17").Value = "'0.776"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.48%  '
$ws.Range("D18").Value = '43.592.77'
$ws.Range("E18").Value = '  +3.05%  '
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("D20").Value = "'69.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("E21").Value = '  +1.43%  '
$ws.Range("E22").Value = '  +15.56%  '
$ws.Range("D23").Value = "'230.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("E24").Value = '  -5.20%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = "'3.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.69%  '
$ws.Range("B27").Value = 'InjectiveProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D27").Value = "'42.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +15.18%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = "'10.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.54%  '
$ws.Range("B29").Value = 'WEMIXToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D29").Value = "'3.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'2.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.63%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = "'2.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.84%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = "'173.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.30%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = "'20.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.33%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = "'0.0876"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.64%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = "'5.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.50%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = "'0.114"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.85%  '
$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").Value = "'0.122"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.13%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = "'4.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.75%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'0.0353"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.48%  '
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").Value = "'13.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.55%  '
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").Value = "'2.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +17.14%  '
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").Value = "'2.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.85%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").Value = "'64.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.22%  '
$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").Value = "'5.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.32%  '
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = "'0.200"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.70%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'100.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = "'0.0979"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = "'8.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.84%  '
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").Value = "'1.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.10%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").Value = "'1.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.92%  '
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").Value = "'0.439"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.23%  '
